$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Aristoteles)
$ws.Range("G4").Value = 104
$ws.Range("H4").Value = 1151
$ws.Range("I4").Value = 1019
$ws.Range("J4").Value = 104
$ws.Range("K4").Value = 224
$ws.Range("L4").Value = 272
$ws.Range("M4").Value = 582
$ws.Range("N4").Value = 653
$ws.Range("O4").Value = 1019
$ws.Range("P4").Value = 1151
$ws.Range("Q4").Value = 1265
$ws.Range("R4").Value = 1313
$ws.Range("V4").Value = -16.81
$ws.Range("W4").Value = 99.58

# Row 18 (MarliRosa)
$ws.Range("G18").Value = 38
$ws.Range("H18").Value = 942
$ws.Range("I18").Value = 834
$ws.Range("J18").Value = 38
$ws.Range("O18").Value = 834
$ws.Range("P18").Value = 942

# Row 20 (Neusa)
$ws.Range("G20").Value = 112
$ws.Range("H20").Value = 1218
$ws.Range("I20").Value = 1102
$ws.Range("J20").Value = 112
$ws.Range("O20").Value = 1102
$ws.Range("P20").Value = 1218

# Row 24 (Diogenes)
$ws.Range("G24").Value = 112
$ws.Range("H24").Value = 1159
$ws.Range("I24").Value = 1030
$ws.Range("J24").Value = 112
$ws.Range("O24").Value = 1030
$ws.Range("P24").Value = 1159

# Row 36 (Dirceujr)
$ws.Range("G36").Value = 173
$ws.Range("H36").Value = 1476
$ws.Range("I36").Value = 1388
$ws.Range("J36").Value = 173
$ws.Range("K36").Value = 259
$ws.Range("L36").Value = 336
$ws.Range("M36").Value = 600
$ws.Range("N36").Value = 674
$ws.Range("O36").Value = 1388
$ws.Range("P36").Value = 1476
$ws.Range("Q36").Value = 1561
$ws.Range("R36").Value = 1638
$ws.Range("V36").Value = -23.03
$ws.Range("W36").Value = 38.71

# Row 45 (Mariarosa)
$ws.Range("G45").Value = 42
$ws.Range("H45").Value = 1076
$ws.Range("J45").Value = 42
$ws.Range("P45").Value = 1076
